$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header (A1) stays "EMP NAME" text-wise; re-assign so the shared-string
# table drops the stale "sunny"/"no"/"yes" entries once they're unused.
$ws.Range("A1").Value = "EMP NAME"

# Rows 2:32 - relabel the employee-name column, shift check in/out dates
# back 24 days (Jan 2025 instead of Jan/Feb 2025), and clear the Offdays
# flag column (no longer "yes"/"no").
for ($r = 2; $r -le 32; $r++) {
    $ws.Cells.Item($r, 1).Value = "Employee Name"

    $checkInDate = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 2).Value2 = $checkInDate - 24

    $checkOutDate = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 4).Value2 = $checkOutDate - 24
}

$ws.Range("F2:F32").ClearContents() | Out-Null

# Match the saved selection state.
$ws.Range("A2:A32").Select() | Out-Null
